# Generate Report for Handoff
# Updates the localization-status report: the file "6104303d-...md" has moved
# out of "Handed back: in sync with en-US" status into "Ready for handoff",
# and the two other files ("ffff5f0b..." and "ffffff7a...") are re-ordered
# ahead of it, now both carrying the "acd3cc98-..." handoff/target/handback
# references that were already in sync.

$wb = $excel.ActiveWorkbook

# ---------- Sheet "Overview" ----------
$ws1 = $wb.Worksheets.Item("Overview")

$ws1.Range("A2").Value = "ffff5f0b0cd7-2e4d-4fac-8347-038a6aff8eeb.md"

$ws1.Range("A3").Value = "ffffff7a5bff45-0784-4e83-b89b-0b0d1af02392.md"

$ws1.Range("A4").Value = "6104303d-37b8-4b92-8309-bfe68b998bc1.md"
$ws1.Range("B4").Value = "Ready for handoff"
$ws1.Range("C4").Value = "Ready for handoff"

# ---------- Sheet "zh-cn" ----------
$ws2 = $wb.Worksheets.Item("zh-cn")

$ws2.Range("A2").Value = "ffff5f0b0cd7-2e4d-4fac-8347-038a6aff8eeb.md"
$ws2.Range("C2").Value = "acd3cc98-32a6-43cb-9bfd-62d79904db49.aab57bf76b38a394b6610a9034d9b6ef5852519a.zh-cn.xlf"
$ws2.Range("D2").Value = "2016-03-09 10:06:00"
$ws2.Range("E2").Value = "acd3cc98-32a6-43cb-9bfd-62d79904db49.md"
$ws2.Range("F2").Value = "acd3cc98-32a6-43cb-9bfd-62d79904db49.aab57bf76b38a394b6610a9034d9b6ef5852519a.zh-cn.xlf"
$ws2.Range("G2").Value = "2016-03-09 10:06:29"

$ws2.Range("A3").Value = "ffffff7a5bff45-0784-4e83-b89b-0b0d1af02392.md"

$ws2.Range("A4").Value = "6104303d-37b8-4b92-8309-bfe68b998bc1.md"
$ws2.Range("B4").Value = "Ready for handoff"
$ws2.Range("C4").Value = "6104303d-37b8-4b92-8309-bfe68b998bc1.507280672b761ae06d2d43713a27199e770b6384.zh-cn.xlf"
$ws2.Range("D4").Value = "2016-03-09 10:08:00"
$ws2.Range("E4").Value = "6104303d-37b8-4b92-8309-bfe68b998bc1.md"
$ws2.Range("F4").Value = "6104303d-37b8-4b92-8309-bfe68b998bc1.507280672b761ae06d2d43713a27199e770b6384.zh-cn.xlf"
$ws2.Range("G4").Value = "2016-03-09 10:07:38"

# ---------- Sheet "de-de" ----------
$ws3 = $wb.Worksheets.Item("de-de")

$ws3.Range("A2").Value = "ffff5f0b0cd7-2e4d-4fac-8347-038a6aff8eeb.md"
$ws3.Range("C2").Value = "acd3cc98-32a6-43cb-9bfd-62d79904db49.aab57bf76b38a394b6610a9034d9b6ef5852519a.de-de.xlf"
$ws3.Range("D2").Value = "2016-03-09 10:06:06"
$ws3.Range("E2").Value = "acd3cc98-32a6-43cb-9bfd-62d79904db49.md"
$ws3.Range("F2").Value = "acd3cc98-32a6-43cb-9bfd-62d79904db49.aab57bf76b38a394b6610a9034d9b6ef5852519a.de-de.xlf"
$ws3.Range("G2").Value = "2016-03-09 10:06:35"

$ws3.Range("A3").Value = "ffffff7a5bff45-0784-4e83-b89b-0b0d1af02392.md"

$ws3.Range("A4").Value = "6104303d-37b8-4b92-8309-bfe68b998bc1.md"
$ws3.Range("B4").Value = "Ready for handoff"
$ws3.Range("C4").Value = "6104303d-37b8-4b92-8309-bfe68b998bc1.507280672b761ae06d2d43713a27199e770b6384.de-de.xlf"
$ws3.Range("D4").Value = "2016-03-09 10:08:07"
$ws3.Range("E4").Value = "6104303d-37b8-4b92-8309-bfe68b998bc1.md"
$ws3.Range("F4").Value = "6104303d-37b8-4b92-8309-bfe68b998bc1.507280672b761ae06d2d43713a27199e770b6384.de-de.xlf"
$ws3.Range("G4").Value = "2016-03-09 10:07:44"

Write-Host "Updated localization-status report for handoff."
